# Automatische test-sync: 2025-06-24 20:02:50
# Appends the newly received "Ruilen van product" mail log entry to the
# "Logs" sheet, extends the conditional formatting ranges to cover the
# new row, and bumps the "Retour / Terugbetaling" tally on the "Dashboard"
# sheet from 4 to 5.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 14 -----------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A14").Value = "Ruilen van product"
$logs.Range("B14").Value = "mailmind.test@zohomail.eu"
$logs.Range("C14").Value = "Kan ik dit product ruilen voor een andere maat?"
$logs.Range("D14").Value = "Retour / Terugbetaling"
$logs.Range("F14").Value = "2025-06-24 20:01:57"
$logs.Range("G14").Value = "Nee"

# --- Extend conditional formatting ranges D2:D13 -> D2:D14 and
#     G2:G13 -> G2:G14 so the new row gets highlighted too -------------
$logs.Range("D2:D13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D14"))
$logs.Range("G2:G13").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G14"))

# --- Dashboard sheet: update the "Retour / Terugbetaling" count -------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 5
